# Dodan opis k projektu
# - shorten several row labels in column A
# - change the report title in A1 to "pridelek"
# - remove the "Deteljno travne mesanice (1 do 5 let)" row (old row 15)
# - select A25 as the last active cell

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title cell
$ws.Range("A1").Value = "pridelek"

# Shorten crop-name labels (rows keyed by their CURRENT row number,
# i.e. before the row-15 deletion below)
$ws.Range("A2").Value  = "Pšenica "
$ws.Range("A6").Value  = "Riž"
$ws.Range("A7").Value  = "Koruza "
$ws.Range("A8").Value  = "Silažna"
$ws.Range("A11").Value = "Repica"
$ws.Range("A13").Value = "Trave"
$ws.Range("A14").Value = "Deteljne"
$ws.Range("A18").Value = "Trajni travniki"
$ws.Range("A19").Value = "Zelje"
$ws.Range("A21").Value = "Jabolka "
$ws.Range("A23").Value = "Breskve"

# Remove the "Deteljno travne mešanice (1 do 5 let)" row entirely
$ws.Rows(15).Delete()

# Leave the selection on A25, matching the saved view state
$ws.Range("A25").Select()
